$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 33.83768468314501
$ws.Range("G2").Value = 0.01117150690685271
$ws.Range("H2").Value = 103.3791977466574
$ws.Range("I2").Value = 0.000405949375518209
$ws.Range("J2").Value = [double]"1.529854878156235e-07"
$ws.Range("K2").Value = 0.001221913840494738
$ws.Range("L2").Value = 0.006277671101953943
$ws.Range("M2").Value = [double]"2.052328410073348e-06"
$ws.Range("N2").Value = 0.01914932138598907
$ws.Range("F3").Value = 0.04777680765572575
$ws.Range("G3").Value = 0.04744115423198638
$ws.Range("H3").Value = 0.04809985782765978
$ws.Range("I3").Value = 0.04617902485631615
$ws.Range("J3").Value = 0.0458557069314974
$ws.Range("K3").Value = 0.04648904698588573
$ws.Range("L3").Value = 0.04780954926117981
$ws.Range("M3").Value = 0.04747402864087715
$ws.Range("N3").Value = 0.04813265109438832
$ws.Range("F4").Value = 33.88546149080073
$ws.Range("G4").Value = 0.0586126611388391
$ws.Range("H4").Value = 103.427297604485
$ws.Range("I4").Value = 0.04658497423183437
$ws.Range("J4").Value = 0.04585585991698522
$ws.Range("K4").Value = 0.04771096082638047
$ws.Range("L4").Value = 0.05408722036313375
$ws.Range("M4").Value = 0.04747608096928723
$ws.Range("N4").Value = 0.0672819724803774
